$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.516.41"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").Value = "1.796.30"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.04"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.25%  "
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.60"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.98%  "
$ws.Range("E9").Value = "  +1.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0694"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.55%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0950"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.30%  "
$ws.Range("D12").Value = "2.055.97"
$ws.Range("E12").Value = "  +0.42%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.824.57"
$ws.Range("E13").Value = "  +1.86%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.11"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.43%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.640"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.92%  "
$ws.Range("D16").Value = "34.497.22"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "247.53"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").Value = "0.0₃0800"
$ws.Range("E20").Value = "  +2.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.23"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.68%  "
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.78%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "163.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("E26").Value = "  +1.52%  "
$ws.Range("E27").Value = "  +1.39%  "
$ws.Range("E28").Value = "  +2.41%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  +0.32%  "
$ws.Range("E31").Value = "  +3.90%  "
$ws.Range("E32").Value = "  +0.60%  "
$ws.Range("E33").Value = "  +7.35%  "
$ws.Range("E34").Value = "  +1.52%  "
$ws.Range("D35").Value = "1.444.64"
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.60"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.16%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.670"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +3.23%  "
$ws.Range("E38").Value = "  +2.76%  "
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "84.34"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.94%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.936"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.57%  "
$ws.Range("E43").Value = "  +2.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.67"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.42%  "
$ws.Range("E45").Value = "  +3.03%  "
$ws.Range("E46").Value = "  +0.70%  "
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("D48").Value = "1.952.06"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "106.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "0.0₆0130"
$ws.Range("E50").Value = "  -4.06%  "
$ws.Range("E51").Value = "  -0.03%  "
